$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the acceptance-criteria text that used to live in B8 ("UI must be
# capable of testing implemented functions") — the row stays, the cell is
# now blank (style retained).
$ws.Range("B8").ClearContents() | Out-Null

# Excel leaves the selection on the cell that was last edited.
$ws.Range("B8").Select() | Out-Null
